$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows right before the current row 208, shifting the
# existing rows 208-244 down to 212-248.
$ws.Rows.Item(208).Resize(4).Insert()

# Values for the 4 freshly-inserted rows (A..R), taken from the new weekly
# price observations added to the dataset.
$newRows = @(
  @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44504, 13, 100112032, "Zapallo italiano", "Bola 8", "Primera", 230, 5000, 6000, 5435, "`$/caja 60 unidades", "Región de Arica y Parinacota", 91, 60, "Hortaliza"),
  @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44504, 13, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 2500, 4000, 5000, 4400, "`$/caja 50 unidades", "Región de Arica y Parinacota", 88, 50, "Hortaliza"),
  @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44504, 13, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 1050, 6000, 7000, 6429, "`$/caja 50 unidades", "Región de O'Higgins", 129, 50, "Hortaliza"),
  @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44504, 13, 100112032, "Zapallo italiano", "Sin especificar", "Segunda", 300, 4000, 5000, 4600, "`$/caja 80 unidades", "Región de O'Higgins", 58, 80, "Hortaliza")
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
  $rowNum = 208 + $i
  $rowVals = $newRows[$i]
  for ($col = 1; $col -le $rowVals.Length; $col++) {
    $ws.Cells.Item($rowNum, $col).Value = $rowVals[$col - 1]
  }
  # Match the date-formatted style used by the rest of column D.
  $ws.Cells.Item($rowNum, 4).NumberFormat = $ws.Cells.Item($rowNum + 4, 4).NumberFormat
}
